$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2506738544474393
$ws.Range("C2").Value = 0.4582210242587601
$ws.Range("J2").Value = 0.01347708894878706
$ws.Range("P2").Value = 0.1752021563342318
$ws.Range("S2").Value = 0.1024258760107817
$ws.Range("B3").Value = 0.01685393258426966
$ws.Range("C3").Value = 0.02247191011235955
$ws.Range("J3").Value = 0.02247191011235955
$ws.Range("P3").Value = 0.7921348314606742
$ws.Range("S3").Value = 0.1460674157303371
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.8444444444444444
$ws.Range("S4").Value = 0.1333333333333333
$ws.Range("B6").Value = 0.09818181818181818
$ws.Range("D6").Value = 0.007272727272727273
$ws.Range("E6").Value = 0.003636363636363636
$ws.Range("F6").Value = 0.08363636363636363
$ws.Range("J6").Value = 0.1490909090909091
$ws.Range("O6").Value = 0.02909090909090909
$ws.Range("Q6").Value = 0.2
$ws.Range("R6").Value = 0.05090909090909091
$ws.Range("S6").Value = 0.3781818181818182
$ws.Range("B7").Value = 0.1218274111675127
$ws.Range("D7").Value = 0.01522842639593909
$ws.Range("F7").Value = 0.03045685279187817
$ws.Range("J7").Value = 0.1522842639593909
$ws.Range("O7").Value = 0.02538071065989848
$ws.Range("Q7").Value = 0.2284263959390863
$ws.Range("R7").Value = 0.02030456852791878
$ws.Range("S7").Value = 0.4060913705583756
$ws.Range("B8").Value = 0.09975669099756691
$ws.Range("D8").Value = 0.0340632603406326
$ws.Range("F8").Value = 0.09002433090024331
$ws.Range("J8").Value = 0.1313868613138686
$ws.Range("O8").Value = 0.03163017031630171
$ws.Range("Q8").Value = 0.1557177615571776
$ws.Range("R8").Value = 0.06082725060827251
$ws.Range("S8").Value = 0.3965936739659368
$ws.Range("B9").Value = 0.1017699115044248
$ws.Range("D9").Value = 0.008849557522123894
$ws.Range("F9").Value = 0.1061946902654867
$ws.Range("J9").Value = 0.1017699115044248
$ws.Range("O9").Value = 0.03097345132743363
$ws.Range("Q9").Value = 0.1902654867256637
$ws.Range("R9").Value = 0.07964601769911504
$ws.Range("S9").Value = 0.3805309734513274
$ws.Range("B10").Value = 0.1134699853587116
$ws.Range("D10").Value = 0.01903367496339678
$ws.Range("F10").Value = 0.08125915080527087
$ws.Range("J10").Value = 0.1295754026354319
$ws.Range("O10").Value = 0.02635431918008785
$ws.Range("Q10").Value = 0.2291361639824304
$ws.Range("R10").Value = 0.05051244509516838
$ws.Range("S10").Value = 0.3506588579795022
$ws.Range("G11").Value = 0.1717791411042945
$ws.Range("J11").Value = 0.09815950920245399
$ws.Range("K11").Value = 0.2208588957055215
$ws.Range("L11").Value = 0.4938650306748466
$ws.Range("S11").Value = 0.01533742331288344
$ws.Range("G12").Value = 0.7716049382716049
$ws.Range("J12").Value = 0.1604938271604938
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.01851851851851852
$ws.Range("S12").Value = 0.03703703703703703
$ws.Range("G13").Value = 0.71875
$ws.Range("J13").Value = 0.21875
$ws.Range("S13").Value = 0.0625
$ws.Range("F15").Value = 0.03501945525291829
$ws.Range("H15").Value = 0.132295719844358
$ws.Range("I15").Value = 0.08949416342412451
$ws.Range("J15").Value = 0.3852140077821012
$ws.Range("K15").Value = 0.07003891050583658
$ws.Range("M15").Value = 0.01167315175097276
$ws.Range("O15").Value = 0.05058365758754864
$ws.Range("S15").Value = 0.2256809338521401
$ws.Range("F16").Value = 0.02542372881355932
$ws.Range("H16").Value = 0.1567796610169492
$ws.Range("I16").Value = 0.0847457627118644
$ws.Range("J16").Value = 0.4576271186440678
$ws.Range("K16").Value = 0.1313559322033898
$ws.Range("M16").Value = 0.00423728813559322
$ws.Range("O16").Value = 0.05932203389830509
$ws.Range("S16").Value = 0.08050847457627118
$ws.Range("F17").Value = 0.01937984496124031
$ws.Range("H17").Value = 0.1511627906976744
$ws.Range("I17").Value = 0.09108527131782945
$ws.Range("J17").Value = 0.4651162790697674
$ws.Range("K17").Value = 0.08914728682170543
$ws.Range("M17").Value = 0.005813953488372093
$ws.Range("N17").Value = 0.001937984496124031
$ws.Range("O17").Value = 0.06395348837209303
$ws.Range("S17").Value = 0.1124031007751938
$ws.Range("F18").Value = 0.0390625
$ws.Range("H18").Value = 0.1640625
$ws.Range("I18").Value = 0.15625
$ws.Range("J18").Value = 0.40625
$ws.Range("K18").Value = 0.0546875
$ws.Range("O18").Value = 0.0703125
$ws.Range("S18").Value = 0.109375
$ws.Range("F19").Value = 0.01511535401750199
$ws.Range("H19").Value = 0.192521877486078
$ws.Range("I19").Value = 0.09148766905330151
$ws.Range("J19").Value = 0.3874303898170247
$ws.Range("K19").Value = 0.1121718377088305
$ws.Range("M19").Value = 0.02147971360381861
$ws.Range("N19").Value = 0.0007955449482895784
$ws.Range("O19").Value = 0.07478122513922036
$ws.Range("S19").Value = 0.1042163882259348
